$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.358412
$ws.Range("B3").Value = 0.1828
$ws.Range("B4").Value = 0.4336
$ws.Range("B5").Value = 2.0676332
